$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2084690553745928
$ws.Range("C2").Value = 0.5602605863192183
$ws.Range("J2").Value = 0.01954397394136808
$ws.Range("P2").Value = 0.1498371335504886
$ws.Range("S2").Value = 0.06188925081433225
# Row 3
$ws.Range("B3").Value = 0.01648351648351648
$ws.Range("C3").Value = 0.04395604395604396
$ws.Range("J3").Value = 0.03296703296703297
$ws.Range("P3").Value = 0.7307692307692307
$ws.Range("S3").Value = 0.1758241758241758
# Row 4
$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("P4").Value = 0.6744186046511628
$ws.Range("S4").Value = 0.2790697674418605
# Row 6
$ws.Range("B6").Value = 0.0947867298578199
$ws.Range("D6").Value = 0.01895734597156398
$ws.Range("E6").Value = 0.004739336492890996
$ws.Range("F6").Value = 0.04265402843601896
$ws.Range("J6").Value = 0.2511848341232227
$ws.Range("O6").Value = 0.04265402843601896
$ws.Range("Q6").Value = 0.1184834123222749
$ws.Range("R6").Value = 0.08056872037914692
$ws.Range("S6").Value = 0.3459715639810427
# Row 7
$ws.Range("B7").Value = 0.09375
$ws.Range("D7").Value = 0.01041666666666667
$ws.Range("E7").Value = 0.005208333333333333
$ws.Range("F7").Value = 0.04166666666666666
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("O7").Value = 0.01041666666666667
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.1041666666666667
$ws.Range("S7").Value = 0.4010416666666667
# Row 8
$ws.Range("B8").Value = 0.0815450643776824
$ws.Range("D8").Value = 0.02145922746781116
$ws.Range("F8").Value = 0.06437768240343347
$ws.Range("J8").Value = 0.09442060085836911
$ws.Range("O8").Value = 0.01716738197424893
$ws.Range("Q8").Value = 0.1502145922746781
$ws.Range("R8").Value = 0.1266094420600858
$ws.Range("S8").Value = 0.444206008583691
# Row 9
$ws.Range("B9").Value = 0.0975609756097561
$ws.Range("D9").Value = 0.01626016260162602
$ws.Range("F9").Value = 0.1138211382113821
$ws.Range("J9").Value = 0.1219512195121951
$ws.Range("Q9").Value = 0.0975609756097561
$ws.Range("R9").Value = 0.1138211382113821
$ws.Range("S9").Value = 0.4390243902439024
# Row 10
$ws.Range("B10").Value = 0.1245901639344262
$ws.Range("D10").Value = 0.02131147540983606
$ws.Range("E10").Value = 0.003278688524590164
$ws.Range("F10").Value = 0.07213114754098361
$ws.Range("J10").Value = 0.1229508196721311
$ws.Range("O10").Value = 0.0180327868852459
$ws.Range("Q10").Value = 0.169672131147541
$ws.Range("R10").Value = 0.1008196721311475
$ws.Range("S10").Value = 0.3672131147540983
# Row 11
$ws.Range("G11").Value = 0.1412103746397695
$ws.Range("J11").Value = 0.1123919308357349
$ws.Range("K11").Value = 0.2161383285302594
$ws.Range("L11").Value = 0.4985590778097982
$ws.Range("S11").Value = 0.03170028818443804
# Row 12
$ws.Range("F12").Value = 0.005681818181818182
$ws.Range("G12").Value = 0.6534090909090909
$ws.Range("J12").Value = 0.25
$ws.Range("K12").Value = 0.01704545454545454
$ws.Range("L12").Value = 0.02840909090909091
$ws.Range("S12").Value = 0.04545454545454546
# Row 13
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2666666666666667
$ws.Range("S13").Value = 0.06666666666666667
# Row 15
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.1504854368932039
$ws.Range("I15").Value = 0.0825242718446602
$ws.Range("J15").Value = 0.2961165048543689
$ws.Range("K15").Value = 0.1116504854368932
$ws.Range("M15").Value = 0.01941747572815534
$ws.Range("O15").Value = 0.06310679611650485
$ws.Range("S15").Value = 0.2572815533980582
# Row 16
$ws.Range("F16").Value = 0.02392344497607655
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.04784688995215311
$ws.Range("J16").Value = 0.430622009569378
$ws.Range("K16").Value = 0.1100478468899522
$ws.Range("M16").Value = 0.01913875598086124
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1339712918660287
# Row 17
$ws.Range("F17").Value = 0.01445086705202312
$ws.Range("H17").Value = 0.1705202312138728
$ws.Range("I17").Value = 0.05202312138728324
$ws.Range("J17").Value = 0.4104046242774567
$ws.Range("K17").Value = 0.1127167630057803
$ws.Range("M17").Value = 0.008670520231213872
$ws.Range("O17").Value = 0.07225433526011561
$ws.Range("S17").Value = 0.1589595375722543
# Row 18
$ws.Range("F18").Value = 0.008695652173913044
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.04347826086956522
$ws.Range("J18").Value = 0.391304347826087
$ws.Range("K18").Value = 0.1130434782608696
$ws.Range("M18").Value = 0.03043478260869565
$ws.Range("N18").Value = 0.004347826086956522
$ws.Range("O18").Value = 0.06956521739130435
$ws.Range("S18").Value = 0.1391304347826087
# Row 19
$ws.Range("F19").Value = 0.01252936570086139
$ws.Range("H19").Value = 0.2317932654659358
$ws.Range("I19").Value = 0.05324980422866092
$ws.Range("J19").Value = 0.355520751761942
$ws.Range("K19").Value = 0.1205951448707909
$ws.Range("M19").Value = 0.02192638997650744
$ws.Range("N19").Value = 0.002349256068911512
$ws.Range("O19").Value = 0.05951448707909162
$ws.Range("S19").Value = 0.1425215348472983

Write-Host "Applied 111 cell updates to ULM_B transition matrix"
